# Ready for remote runs
# Convert the dynamic-array "lookup" formulas (FILTER/TRANSPOSE/LET) that live on the
# validation-list helper sheets into their plain literal values (values-only paste),
# then restore each sheet's cursor position and re-select the previously active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) mother_vessels — A2 was an array FILTER formula -> literal "V4"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("mother_vessels")
$ws.Range("A2").Value2 = "V4"
$ws.Range("A3").Select()

# ---------------------------------------------------------------------------
# 2) locations — A2 was a LET/VSTACK array formula -> literal "B1"; A3 -> "V4"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("locations")
$ws.Range("A2").Value2 = "B1"
$ws.Range("A3").Value2 = "V4"
$ws.Range("E2").Select()

# ---------------------------------------------------------------------------
# 3) docking — no formulas, just a cursor move
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("docking")
$ws.Range("E1").Select()

# ---------------------------------------------------------------------------
# 4) tasks — no formulas, just a cursor move
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("tasks")
$ws.Range("A5").Select()

# ---------------------------------------------------------------------------
# 5) spare_parts — no formulas, just a cursor move
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("spare_parts")
$ws.Range("A2").Select()

# ---------------------------------------------------------------------------
# 6) task_compatibility — B1:F1 TRANSPOSE/FILTER -> literals V1..V5;
#    A2:A5 FILTER -> literals M1..M4
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("task_compatibility")
$ws.Range("B1").Value2 = "V1"
$ws.Range("C1").Value2 = "V2"
$ws.Range("D1").Value2 = "V3"
$ws.Range("E1").Value2 = "V4"
$ws.Range("F1").Value2 = "V5"
$ws.Range("A2").Value2 = "M1"
$ws.Range("A3").Value2 = "M2"
$ws.Range("A4").Value2 = "M3"
$ws.Range("A5").Value2 = "M4"
$ws.Range("F16").Select()

# ---------------------------------------------------------------------------
# 7) capacity_base_vessels — B1:F1 TRANSPOSE/FILTER -> literals V1..V5;
#    A2 FILTER -> literal "B1"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("capacity_base_vessels")
$ws.Range("B1").Value2 = "V1"
$ws.Range("C1").Value2 = "V2"
$ws.Range("D1").Value2 = "V3"
$ws.Range("E1").Value2 = "V4"
$ws.Range("F1").Value2 = "V5"
$ws.Range("A2").Value2 = "B1"
$ws.Range("F11").Select()

# ---------------------------------------------------------------------------
# 8) holding_costs — B1:C1 TRANSPOSE/FILTER -> literals B1,V4; A2 FILTER -> "S1"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("holding_costs")
$ws.Range("B1").Value2 = "B1"
$ws.Range("C1").Value2 = "V4"
$ws.Range("A2").Value2 = "S1"
$ws.Range("B2").Select()

# ---------------------------------------------------------------------------
# 9) spare_parts_required — B1:E1 TRANSPOSE/FILTER -> literals M1..M4;
#    A2 FILTER -> "S1"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("spare_parts_required")
$ws.Range("B1").Value2 = "M1"
$ws.Range("C1").Value2 = "M2"
$ws.Range("D1").Value2 = "M3"
$ws.Range("E1").Value2 = "M4"
$ws.Range("A2").Value2 = "S1"
$ws.Range("B2").Select()

# ---------------------------------------------------------------------------
# 10) max_capacity — B1:C1 TRANSPOSE/FILTER -> literals B1,V4; A2 FILTER -> "S1"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("max_capacity")
$ws.Range("B1").Value2 = "B1"
$ws.Range("C1").Value2 = "V4"
$ws.Range("A2").Value2 = "S1"
$ws.Range("D9").Select()

# ---------------------------------------------------------------------------
# 11) reorder_level — B1:C1 TRANSPOSE/FILTER -> literals B1,V4; A2 FILTER -> "S1";
#     the stray spilled D2 cell disappears along with the formula's spill range
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("reorder_level")
$ws.Range("B1").Value2 = "B1"
$ws.Range("C1").Value2 = "V4"
$ws.Range("A2").Value2 = "S1"
$ws.Range("D2").ClearContents()
$ws.Range("H18").Select()

# ---------------------------------------------------------------------------
# 12) Re-point the active sheet from "vessels" to "bases" and restore each of
#     their cursor positions.
# ---------------------------------------------------------------------------
$wsVessels = $wb.Worksheets.Item("vessels")
$wsVessels.Range("E1").Select()

$wsBases = $wb.Worksheets.Item("bases")
$wsBases.Activate()
$wsBases.Range("A2").Select()
